# Applies updated vm_pu results (case with 380 kV done) to rows 2-25,
# columns B-F and I-N, leaving row 1 header, column G (=1) and column A
# (bus index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.019999999999999, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02)
$Cvals = @(1.029642872460143, 1.030600790180702, 1.031220591702128, 1.03148114812495, 1.031524896226788, 1.031224073301386, 1.029966612737224, 1.027750521542847, 1.0262729122298, 1.025633033853964, 1.025395344308621, 1.025446329996064, 1.025613386574336, 1.025716314271446, 1.026315377451511, 1.026691136427108, 1.026910304168546, 1.026985033687364, 1.026650821681794, 1.025564192863277, 1.024880926711031, 1.02524314488055, 1.026669038198494, 1.028323470419253)
$Dvals = @(1.032833735723286, 1.033530064845228, 1.033980588234684, 1.034169976335252, 1.034201774718837, 1.033983118895173, 1.033069073125948, 1.031458056836278, 1.030383831593888, 1.02991863316344, 1.029745830194165, 1.029782897356213, 1.029904349369028, 1.029979178967463, 1.030414704209387, 1.030687883977194, 1.030847220100212, 1.03090154874317, 1.030658574890311, 1.029868584985271, 1.029371843347434, 1.029635179550995, 1.030671818420488, 1.031874582496238)
$Evals = @(1.033278361253346, 1.034184041280296, 1.034770391529595, 1.035016968440821, 1.035058374200582, 1.034773686006711, 1.033584374851946, 1.031491069543372, 1.030097155331567, 1.029493959028182, 1.029269961628853, 1.029318007274787, 1.029475442180478, 1.02957245051157, 1.030137195426548, 1.030491545995715, 1.030698269147851, 1.030768762590403, 1.030453523763129, 1.029429079964782, 1.028785298658538, 1.029126548362431, 1.030470704260557, 1.032031953136657)
$Fvals = @(1.02821394495719, 1.029804137700556, 1.030833035577246, 1.031265574650826, 1.031338199462438, 1.030838815216339, 1.028751374447106, 1.025072305920073, 1.022618755033256, 1.02155606402449, 1.021161283372419, 1.021245967528669, 1.021523432389528, 1.021694380834316, 1.022689275708708, 1.023313265320788, 1.023677200638798, 1.023801288757038, 1.023246320008924, 1.021441727335197, 1.020306816265647, 1.020908484034467, 1.023276569796414, 1.026023555623283)
$Ivals = @(1.034117349114099, 1.034337180944157, 1.034478224356251, 1.034537231101289, 1.034547121709661, 1.034479013938576, 1.034191891275137, 1.033676735872699, 1.03332710766339, 1.033174245884125, 1.033117245196889, 1.033129482040815, 1.033169538698765, 1.033194189638051, 1.033337221607182, 1.033426548031764, 1.033478508752693, 1.033496201954834, 1.033416978826829, 1.033157749113299, 1.03299348232272, 1.033080684502836, 1.033421303176914, 1.033811007144842)
$Jvals = @(1.034788514523573, 1.035387507612278, 1.03577446382193, 1.035936988991353, 1.035964268819649, 1.035776636085191, 1.034991078082914, 1.033601973107908, 1.032672627514841, 1.032269429872461, 1.032119545981737, 1.032151701946701, 1.032257042839704, 1.032321931147617, 1.032699369842788, 1.032935916741157, 1.033073814844811, 1.033120821707079, 1.032910545330508, 1.032226025842237, 1.03179495641721, 1.032023539532474, 1.032922009812588, 1.033961666357509)
$Kvals = @(1.035637825768086, 1.036143310324102, 1.036469670567916, 1.036606699496834, 1.036629697115359, 1.036471502234697, 1.03580880611768, 1.034635522553382, 1.033849616819133, 1.033508428742654, 1.033381563088782, 1.033408782232624, 1.033497944706438, 1.033552862993177, 1.033872241675143, 1.034072342256138, 1.034188972091595, 1.034228725362316, 1.034050882196744, 1.033471692258666, 1.033106761764752, 1.033300291435798, 1.034060579337735, 1.034939499569091)
$Lvals = @(1.03608117050749, 1.036795539243833, 1.037257469763535, 1.037451590297705, 1.037484179590283, 1.037260063906454, 1.036322659790857, 1.034668426530121, 1.033563977500297, 1.033085350326705, 1.032907507443795, 1.032945658040984, 1.033070650992996, 1.033147655389768, 1.033595734092025, 1.03387669645575, 1.034040539108596, 1.034096398793654, 1.033846555779701, 1.033033845327486, 1.032522517942806, 1.032793614988798, 1.03386017517086, 1.035096371598089)
$Mvals = @(1.031031414052521, 1.032427385176808, 1.033330137610621, 1.03370953157025, 1.033773226400503, 1.033335207569347, 1.031503302495856, 1.028270948222689, 1.026112854529444, 1.025177552765333, 1.024830009105135, 1.024904564418704, 1.025148827402737, 1.025299308364127, 1.026174909383365, 1.026723923679199, 1.027044074451546, 1.027153224201925, 1.026665027954808, 1.025076901703322, 1.024077623226097, 1.024607433680659, 1.026691640643273, 1.02910712983142)
$Nvals = @(1.03625803290611, 1.036857876633695, 1.037245382365549, 1.037408138339354, 1.037435456908135, 1.037247557713673, 1.036460884128927, 1.035069806465525, 1.034139141095118, 1.033735370865839, 1.033585274122798, 1.03361747575292, 1.033722966242074, 1.033787946698827, 1.034165921400239, 1.034402804222335, 1.034540898157114, 1.034587971774518, 1.034377396781373, 1.033691905196913, 1.033260223603859, 1.033489131333365, 1.034388877544332, 1.035430010520798)

$colMap = [ordered]@{
    "B" = $Bvals
    "C" = $Cvals
    "D" = $Dvals
    "E" = $Evals
    "F" = $Fvals
    "I" = $Ivals
    "J" = $Jvals
    "K" = $Kvals
    "L" = $Lvals
    "M" = $Mvals
    "N" = $Nvals
}

foreach ($col in $colMap.Keys) {
    $vals = $colMap[$col]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $row = $i + 2
        $ws.Range("$col$row").Value = $vals[$i]
    }
}
